# AGA206 Assessment 2 Checklist - "new level and pause menu"
#
# The author ticked three checkboxes on the checklist:
#   - J15 (Task 4  - "Pause & Restart Screen")
#   - J23 (Task 10 - "Level Select Menu & Multiple Levels")
#   - J24 (Task 11 - "Speed Run Mode")
#
# These are legacy Forms checkboxes whose `LinkedCell` (fmlaLink) points at
# J15 / J23 / J24, so ticking the box == writing TRUE into the linked cell.
# Writing the cells directly reproduces every dependent-formula ripple
# (F15/F23/F24 "Done" text, K23/K24 per-row points, D7/D8/D9 summary counts,
# K39 total) exactly the way Excel's recalculation engine would.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

$ws.Range("J15").Value = $true
$ws.Range("J23").Value = $true
$ws.Range("J24").Value = $true

# Also flip the on-sheet checkbox shapes themselves so their state/appearance
# stays consistent with the linked cells (best effort - the underlying
# control metadata is otherwise untouched by just writing the cell).
$checkboxNames = @("Check Box 14", "Check Box 36", "Check Box 38")
foreach ($name in $checkboxNames) {
    try {
        $shp = $ws.Shapes.Item($name)
        $shp.ControlFormat.Value = 1
    } catch {
        # ignore - shape lookup/version differences shouldn't abort the edit
    }
}

# The author had scrolled down to/selected the newly-updated rows before
# saving (sheetView selection moved from G47 to G25, with the window
# scrolled so row 25 is at the top).
$ws.Activate() | Out-Null
$target = $ws.Range("G25")
$target.Select() | Out-Null
$excel.Goto($target, $true) | Out-Null
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1

$wb.Application.Calculate()
